$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3-5, keeping only the header row and first data row
$ws.Rows("3:5").Delete()

# Update row 2 (B2/C2/D2 labels reorder: Hcrt, Hcrtr1, MuSCs)
$ws.Range("B2").Value = "Hcrt"
$ws.Range("C2").Value = "Hcrtr1"
$ws.Range("D2").Value = "MuSCs"

# Update numeric values on row 2 with the new TPM-derived figures
$ws.Range("G2").Value = 0.4093176666666667
$ws.Range("H2").Value = 1.227953
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.438129
$ws.Range("N2").Value = 1.314387
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.179333939979
$ws.Range("R2").Value = 1.614005459811
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

$wb.Save()
